# eirini updated catalog with 98% coverage
#
# This script reproduces (on Sheet1) the effect-size / standard-error
# rework of the meta-analysis catalogue:
#   - header relabelled: "WMD" -> "effect", "SE" -> "se"
#   - new "se" column (J) computed from the existing CI bounds:
#       se = (CIH - CIL) / 1.96^2
#   - a handful of previously-blank "id" cells (D column, the rows that
#     continue a multi-arm study) are filled in
#   - a new annotation cell explaining the "effect" column is "WMD"
#   - selection moved to where the editor left off

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header relabel ---------------------------------------------------
$ws.Range("G1").Value = "effect"
$ws.Range("J1").Value = "se"

# --- Fill in missing continuation "id" values in column D -------------
$ws.Range("D4").Value = 2
$ws.Range("D6").Value = 3
$ws.Range("D10").Value = 6
$ws.Range("D15").Value = 10
$ws.Range("D18").Value = 12

# --- Compute se = (CIH - CIL) / 1.96^2 for every data row (2-24) ------
# J2/J3 are entered individually, then J4:J24 filled down as one block -
# this mirrors how the column was actually authored (and how Excel
# groups the fill into a single shared formula).
$ws.Range("J2").Formula = "=(I2-H2)/(1.96^2)"
$ws.Range("J3").Formula = "=(I3-H3)/(1.96^2)"
$ws.Range("J4:J24").Formula = "=(I4-H4)/(1.96^2)"

# --- New annotation cell explaining the effect measure -----------------
$ws.Range("N11").Value = "effect=WMD"
$ws.Range("N11").Style = "Normal"

# --- Restore the editor's last selection --------------------------------
$ws.Range("D28").Select()
